$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update worker identification (Documento / Nombre / Periodo)
$ws.Range("C16").Value = "45552629"
$ws.Range("D16").Value = "DARLYS PATERNINA CANTILLO"
$ws.Range("E16").Value = "2508"

# Update "Valor Mora" amounts (header summary + detail row)
$ws.Range("E11").Value = 47450
$ws.Range("F16").Value = 47450

# Column D width was recalculated by Excel's best-fit after the name changed
$ws.Columns.Item(4).ColumnWidth = 27.3
